$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 50; $r++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -and $val -ne "N/A") {
            $parts = $val -split ":"
            $newParts = @()
            foreach ($p in $parts) {
                $newParts += "0x$p"
            }
            $newVal = $newParts -join ":"
            $cell.Value = $newVal
        }
    }
}
